$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)
$sh = $s.Shapes.Item(4)
$nbsp = [char]0x00A0

$tf = $sh.TextFrame
$tr = $tf.TextRange
$tr.Text = "                                  +------------+`r                                  | Controller |`r                                  +------------+`r   PLM Mode                           /    \      Timestamp Label/SRV6 EB         `r   Loopback or Enhanced Mode         /      \       Timestamp Offse`r   Timestamp Label/SRv6 EB          /        \      Timestamp Format`r     Timestamp Format              /          \`r   Missed Packet Count (N)        /            \`r   Delay Threshold/Count (T/M)   /              \`r   Packet Loss Threshold (XofY) /                \`r                               v                  v`r                           +-------+          +-------+`r                           |       |          |       |`r                           |   R1  |==========|   R3  |`r                           |       |          |       |`r                           +-------+          +-------+`r" + $nbsp + "                        `r                        Session-Sender      Session-Reflector`r" + $nbsp + "  `r" + $nbsp + "                        Figure: Example Provisioning Model"

# Restore shape size: width grows to 516pt (6553200 EMU), height stays 273.84842519685037pt (3477875 EMU)
$sh.Width = 516.0
$sh.Height = 273.84843519685035

# Re-apply blue color (srgbClr 0070C0 -> COM BGR RGB 0xC07000) to the two "Timestamp Label" callouts
$full = $tr.Text
$i1 = $full.IndexOf("Timestamp Label/SRV6 EB")
$tr.Characters($i1 + 1, 23).Font.Color.RGB = 0xC07000

$i2 = $full.IndexOf("Timestamp Label/SRv6 EB")
$tr.Characters($i2 + 1, 33).Font.Color.RGB = 0xC07000

Write-Host "Edit applied successfully"
